# Updated cryptos list on Wed May 31 19:45:12 UTC 2023 with GitHub Actions
#
# Applies the per-cell Price (D) / Volume(1h) (E) refresh, plus the
# WrappedEther/TRON row swap in rows 12-13, as captured by the source diff.
# Price values are written as literal text (leading "'" + reset to the
# "Normal" style) so Excel does not silently reinterpret strings such as
# "1.000" or "0.07603" as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    # Force literal text even for values that look numeric (e.g. "1.000"),
    # then strip the resulting style so no extra formatting is introduced.
    $ws.Range($Address).Value = "'" + $Text
    $ws.Range($Address).Style = 'Normal'
}

# --- Row 2 (Bitcoin) ---
$ws.Range('D2').Value = '27.068.05'
$ws.Range('E2').Value = '  -2.72%  '

# --- Row 3 (Ethereum) ---
$ws.Range('D3').Value = '1.867.18'
$ws.Range('E3').Value = '  -2.09%  '

# --- Row 4 (TetherUSD) ---
Set-TextValue 'D4' '1.000'
$ws.Range('E4').Value = '  +0.07%  '

# --- Row 5 (BNB) ---
Set-TextValue 'D5' '306.02'
$ws.Range('E5').Value = '  -2.17%  '

# --- Row 6 (USDC) ---
Set-TextValue 'D6' '0.9997'
$ws.Range('E6').Value = '  +0.01%  '

# --- Row 7 (XRP) ---
Set-TextValue 'D7' '0.5132'
$ws.Range('E7').Value = '  -1.90%  '

# --- Row 8 (Cardano) ---
Set-TextValue 'D8' '0.3757'
$ws.Range('E8').Value = '  -0.82%  '

# --- Row 9 (Dogecoin) ---
Set-TextValue 'D9' '0.07161'
$ws.Range('E9').Value = '  -1.00%  '

# --- Row 10 (Polygon) ---
Set-TextValue 'D10' '0.8898'
$ws.Range('E10').Value = '  -2.02%  '

# --- Row 11 (Solana) ---
$ws.Range('E11').Value = '  -2.96%  '

# --- Row 12: was WrappedEther, now TRON ---
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D12' '0.07603'
$ws.Range('E12').Value = '  -0.46%  '

# --- Row 13: was TRON, now WrappedEther ---
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.869.94'
$ws.Range('E13').Value = '  -2.14%  '

# --- Row 14 (Polkadot) ---
Set-TextValue 'D14' '5.298'
$ws.Range('E14').Value = '  -2.83%  '

# --- Row 15 (Litecoin) ---
Set-TextValue 'D15' '89.67'
$ws.Range('E15').Value = '  -2.82%  '

# --- Row 16 (BinanceUSD) ---
Set-TextValue 'D16' '1.001'
$ws.Range('E16').Value = '  +0.09%  '

# --- Row 17 (ShibaInu) ---
Set-TextValue 'D17' '0.000008472'
$ws.Range('E17').Value = '  -2.65%  '

# --- Row 18 (Avalanche) ---
$ws.Range('E18').Value = '  -3.21%  '

# --- Row 19 (Dai) ---
Set-TextValue 'D19' '1.000'
$ws.Range('E19').Value = '  +0.09%  '

# --- Row 20 (WrappedBTC) ---
$ws.Range('D20').Value = '27.092.17'
$ws.Range('E20').Value = '  -2.71%  '

# --- Row 21 (Uniswap) ---
Set-TextValue 'D21' '5.033'
$ws.Range('E21').Value = '  -2.27%  '

# --- Row 22 (WrappedliquidstakedEther2.0) ---
$ws.Range('D22').Value = '2.085.31'
$ws.Range('E22').Value = '  -3.50%  '

# --- Row 23 (Cosmos) ---
Set-TextValue 'D23' '10.51'

# --- Row 24 (Chainlink) ---
Set-TextValue 'D24' '6.459'
$ws.Range('E24').Value = '  -2.53%  '

# --- Row 25 (Toncoin) ---
Set-TextValue 'D25' '1.842'
$ws.Range('E25').Value = '  -1.35%  '

# --- Row 26 (Monero) ---
Set-TextValue 'D26' '147.56'
$ws.Range('E26').Value = '  -4.16%  '

# --- Row 27 (EthereumClassic) ---
$ws.Range('E27').Value = '  -1.88%  '

# --- Row 28 (LidoDAOToken) ---
Set-TextValue 'D28' '2.111'
$ws.Range('E28').Value = '  -2.79%  '

# --- Row 29 (BitcoinCash) ---
Set-TextValue 'D29' '112.68'
$ws.Range('E29').Value = '  -1.60%  '

# --- Row 30 (InternetComputer(DFINITY)) ---
Set-TextValue 'D30' '4.657'
$ws.Range('E30').Value = '  -4.01%  '

# --- Row 31 (Filecoin) ---
$ws.Range('E31').Value = '  -3.40%  '

# --- Row 32 (Stellar) ---
Set-TextValue 'D32' '0.09114'
$ws.Range('E32').Value = '  +1.04%  '

# --- Row 33 (Hedera) ---
Set-TextValue 'D33' '0.05127'
$ws.Range('E33').Value = '  -2.94%  '

# --- Row 34 (HuobiToken) ---
Set-TextValue 'D34' '3.068'
$ws.Range('E34').Value = '  -3.40%  '

# --- Row 35 (ARBITRUM) ---
$ws.Range('E35').Value = '  -5.91%  '

# --- Row 36 (ImmutableX) ---
Set-TextValue 'D36' '0.7256'
$ws.Range('E36').Value = '  -7.09%  '

# --- Row 37 (VeChain) ---
$ws.Range('E37').Value = '  -2.52%  '

# --- Row 38 (RenderToken) ---
Set-TextValue 'D38' '2.504'
$ws.Range('E38').Value = '  -4.46%  '

# --- Row 39 (MXToken) ---
Set-TextValue 'D39' '3.037'
$ws.Range('E39').Value = '  -1.27%  '

# --- Row 40 (TrustWalletToken) ---
Set-TextValue 'D40' '1.074'
$ws.Range('E40').Value = '  -1.79%  '

# --- Row 41 (TheSandbox) ---
Set-TextValue 'D41' '0.5338'
$ws.Range('E41').Value = '  -3.84%  '

# --- Row 42 (FraxShare) ---
Set-TextValue 'D42' '6.566'
$ws.Range('E42').Value = '  -1.96%  '

# --- Row 43 (Quant) ---
Set-TextValue 'D43' '115.88'
$ws.Range('E43').Value = '  +0.97%  '

# --- Row 44 (Aptos) ---
$ws.Range('E44').Value = '  -3.31%  '

# --- Row 45 (Algorand) ---
Set-TextValue 'D45' '0.1466'
$ws.Range('E45').Value = '  -3.04%  '

# --- Row 46 (Decentraland) ---
$ws.Range('E46').Value = '  -3.50%  '

# --- Row 47 (PaxDollar) ---
Set-TextValue 'D47' '0.9995'
$ws.Range('E47').Value = '  +0.02%  '

# --- Row 48 (EnergySwap) ---
Set-TextValue 'D48' '9.982'
$ws.Range('E48').Value = '  -4.47%  '

# --- Row 49 (NEARProtocol) ---
$ws.Range('E49').Value = '  -2.91%  '

# --- Row 50 (Elrond) ---
Set-TextValue 'D50' '36.54'
$ws.Range('E50').Value = '  -1.28%  '

# --- Row 51 (Aave) ---
Set-TextValue 'D51' '63.84'
$ws.Range('E51').Value = '  -4.61%  '
